# PRJ0019297 CVAS Time Tracking Changes+ changes after SB Referesh
#
# The "Users" sheet's single test-data row referenced a staff member
# named "Coartney Williams". Update it to "Coartney Trone", bold the
# header cell (matching the header styling already used on the other
# sheets, e.g. Project_Title's A1/B1), set the page to portrait, and
# make "Users" the active/selected sheet (it was "Project_Title").

$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")

# Rename the Global Search User test value.
$wsUsers.Range("A2").Value = "Coartney Trone"

# Bold the header cell, consistent with the other sheets' header rows.
$wsUsers.Range("A1").Font.Bold = $true

# Give this sheet an explicit (portrait) page setup.
$wsUsers.PageSetup.Orientation = 1

# Make "Users" the active sheet / active tab, cursor back at A1.
[void]$wsUsers.Activate()
[void]$wsUsers.Range("A1").Select()
